$wb = $excel.ActiveWorkbook

# --- Overview sheet: widen columns E and F, and roll up the new status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.9777050018311
$wsOverview.Range("F1").ColumnWidth = 29.9777050018311
$wsOverview.Range("E2:F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: report generated for handback ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").ColumnWidth = 29.9777050018311
$wsZh.Range("P1").ColumnWidth = 13.7470531463623
$wsZh.Range("C2:C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K2:K3").Value = "2016-11-09 07:06:59"
$wsZh.Range("P3").Value = ""

# --- de-de sheet: report generated for handback ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").ColumnWidth = 29.9777050018311
$wsDe.Range("P1").ColumnWidth = 13.7470531463623
$wsDe.Range("C2:C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K2:K3").Value = "2016-11-09 07:07:17"
$wsDe.Range("P3").Value = ""
